$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("units new statlines")
$ws.Range("C24").Value = 1
